$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" column header in H1, matching the formatting of the
# existing header cells (bold font, border, centered alignment - same as G1)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Add the new Save column data values (plain numeric cells, no special style)
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
